# Backup QR Scanner data - 2025-11-19T07:05:49.110Z - Cache Bust: 1763535949110
#
# 1) Rename the sheet from "Session" to "Parasitology_SGD_POS".
# 2) Append 4 new scanner-log rows (52-55) with the same column layout as
#    the existing data (Student ID, Subject, Log Date, Log Time, Type, User).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the worksheet -------------------------------------------------
$ws.Name = "Parasitology_SGD_POS"

# --- 2. Append the new log rows ---------------------------------------------
# Column A holds numeric-looking IDs that must stay stored as text (matching
# every other row on the sheet), so force a text number format on the new
# cells in column A before writing the values.
$ws.Range("A52:A55").NumberFormat = "@"

$newRows = @(
    @("244632", "Parasitology SGD/POS", "19/11/2025", "09:02:26", "Scan",   "Alshimaa_khaled@med.asu.edu.eg"),
    @("244631", "Parasitology SGD/POS", "19/11/2025", "09:02:32", "Scan",   "Alshimaa_khaled@med.asu.edu.eg"),
    @("244571", "Parasitology SGD/POS", "19/11/2025", "09:03:01", "Manual", "Alshimaa_khaled@med.asu.edu.eg"),
    @("244499", "Parasitology SGD/POS", "19/11/2025", "09:03:10", "Manual", "Alshimaa_khaled@med.asu.edu.eg")
)

$r = 52
foreach ($rowValues in $newRows) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
    $r++
}
